{"js": "// Update the date line and all 25 \"two-digit \u00d7 two-digit\" equations in the\n// practice table. Each old value is unique in the document, so a simple\n// search + replace (preserving the existing run formatting) is safe for\n// every entry, including the heading date.\nconst replacements = [\n  [\"2025-05-26 Monday\", \"2025-05-27 Tuesday\"],\n  [\"21\u00d782=\", \"72\u00d770=\"],\n  [\"22\u00d721=\", \"55\u00d738=\"],\n  [\"82\u00d799=\", \"67\u00d746=\"],\n  [\"65\u00d731=\", \"69\u00d763=\"],\n  [\"41\u00d771=\", \"64\u00d787=\"],\n  [\"11\u00d744=\", \"29\u00d751=\"],\n  [\"58\u00d748=\", \"53\u00d762=\"],\n  [\"52\u00d742=\", \"72\u00d778=\"],\n  [\"47\u00d780=\", \"94\u00d774=\"],\n  [\"47\u00d715=\", \"53\u00d725=\"],\n  [\"12\u00d762=\", \"42\u00d787=\"],\n  [\"96\u00d775=\", \"84\u00d735=\"],\n  [\"50\u00d781=\", \"51\u00d797=\"],\n  [\"76\u00d757=\", \"71\u00d745=\"],\n  [\"98\u00d794=\", \"86\u00d731=\"],\n  [\"72\u00d761=\", \"66\u00d738=\"],\n  [\"71\u00d794=\", \"69\u00d771=\"],\n  [\"46\u00d719=\", \"82\u00d793=\"],\n  [\"47\u00d729=\", \"23\u00d722=\"],\n  [\"50\u00d774=\", \"59\u00d786=\"],\n  [\"24\u00d728=\", \"97\u00d756=\"],\n  [\"49\u00d756=\", \"96\u00d768=\"],\n  [\"21\u00d766=\", \"26\u00d723=\"],\n  [\"66\u00d740=\", \"42\u00d739=\"],\n  [\"48\u00d734=\", \"95\u00d799=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each old value is unique in the document, so Find/Replace (wdReplaceOne)\n# targeted at the full document Range is safe for every entry, including\n# the heading date.\n$replacements = @(\n    @{ Old = \"2025-05-26 Monday\"; New = \"2025-05-27 Tuesday\" },\n    @{ Old = \"21\u00d782=\"; New = \"72\u00d770=\" },\n    @{ Old = \"22\u00d721=\"; New = \"55\u00d738=\" },\n    @{ Old = \"82\u00d799=\"; New = \"67\u00d746=\" },\n    @{ Old = \"65\u00d731=\"; New = \"69\u00d763=\" },\n    @{ Old = \"41\u00d771=\"; New = \"64\u00d787=\" },\n    @{ Old = \"11\u00d744=\"; New = \"29\u00d751=\" },\n    @{ Old = \"58\u00d748=\"; New = \"53\u00d762=\" },\n    @{ Old = \"52\u00d742=\"; New = \"72\u00d778=\" },\n    @{ Old = \"47\u00d780=\"; New = \"94\u00d774=\" },\n    @{ Old = \"47\u00d715=\"; New = \"53\u00d725=\" },\n    @{ Old = \"12\u00d762=\"; New = \"42\u00d787=\" },\n    @{ Old = \"96\u00d775=\"; New = \"84\u00d735=\" },\n    @{ Old = \"50\u00d781=\"; New = \"51\u00d797=\" },\n    @{ Old = \"76\u00d757=\"; New = \"71\u00d745=\" },\n    @{ Old = \"98\u00d794=\"; New = \"86\u00d731=\" },\n    @{ Old = \"72\u00d761=\"; New = \"66\u00d738=\" },\n    @{ Old = \"71\u00d794=\"; New = \"69\u00d771=\" },\n    @{ Old = \"46\u00d719=\"; New = \"82\u00d793=\" },\n    @{ Old = \"47\u00d729=\"; New = \"23\u00d722=\" },\n    @{ Old = \"50\u00d774=\"; New = \"59\u00d786=\" },\n    @{ Old = \"24\u00d728=\"; New = \"97\u00d756=\" },\n    @{ Old = \"49\u00d756=\"; New = \"96\u00d768=\" },\n    @{ Old = \"21\u00d766=\"; New = \"26\u00d723=\" },\n    @{ Old = \"66\u00d740=\"; New = \"42\u00d739=\" },\n    @{ Old = \"48\u00d734=\"; New = \"95\u00d799=\" },\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
